# Form 16 2017-18 completion: add ROP column (T) to Sheet2
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column T
$ws.Range("T1").Value = "ROP"

# ROP values for each data row (row 93 intentionally has no ROP value)
$tValues = @{
    2 = 0
    3 = 0
    4 = 112
    5 = 0
    6 = 0
    7 = 0
    8 = 0
    9 = 0
    10 = 0
    11 = 0
    12 = 26425
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 0
    21 = 0
    22 = 0
    23 = 30993
    24 = 0
    25 = 0
    26 = 0
    27 = 0
    28 = 112
    29 = 112
    30 = 112
    31 = 112
    32 = 112
    33 = 112
    34 = 112
    35 = 112
    36 = 3892
    37 = 112
    38 = 28912
    39 = 3892
    40 = 0
    41 = 112
    42 = 644
    43 = 17991
    44 = 112
    45 = 112
    46 = 13696
    47 = 1204
    48 = 112
    49 = 112
    50 = 644
    51 = 112
    52 = 112
    53 = 112
    54 = 0
    55 = 112
    56 = 112
    57 = 112
    58 = 112
    59 = 112
    60 = 112
    61 = 112
    62 = 112
    63 = 112
    64 = 112
    65 = 0
    66 = 112
    67 = 112
    68 = 112
    69 = 11092
    70 = 3892
    71 = 112
    72 = 112
    73 = 4279
    74 = 0
    75 = 0
    76 = 0
    77 = 2602
    78 = 8036
    79 = 0
    80 = 0
    81 = 0
    82 = 0
    83 = 0
    84 = 0
    85 = 0
    86 = 0
    87 = 0
    88 = 0
    89 = 0
    90 = 0
    91 = 0
    92 = 0
    94 = 0
}

foreach ($r in $tValues.Keys) {
    $ws.Cells.Item($r, 20).Value = $tValues[$r]
}

# Match the selection state recorded in the saved workbook
$ws.Range("T1:T1048576").Select()
